$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 452 (old rows 452:465 shift down to 455:468)
$ws.Range("A452:T454").EntireRow.Insert()

# --- Fill the 3 new rows (452-454) with new data ---

# Row 452
$ws.Cells.Item(452,1).Value = 11
$ws.Cells.Item(452,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(452,3).Value = "Bíobío"
$ws.Cells.Item(452,4).Value = 44939
$ws.Cells.Item(452,5).Value = 8
$ws.Cells.Item(452,6).Value = "Fruta"
$ws.Cells.Item(452,7).Value = 100101
$ws.Cells.Item(452,8).Value = "Berries"
$ws.Cells.Item(452,9).Value = 100112025
$ws.Cells.Item(452,10).Value = "Frutilla"
$ws.Cells.Item(452,11).Value = "Sin especificar"
$ws.Cells.Item(452,12).Value = "Especial"
$ws.Cells.Item(452,13).Value = 200
$ws.Cells.Item(452,14).Value = 7500
$ws.Cells.Item(452,15).Value = 7500
$ws.Cells.Item(452,16).Value = 7500
$ws.Cells.Item(452,17).Value = "$/caja 7 kilos"
$ws.Cells.Item(452,18).Value = "Región del Maule"
$ws.Cells.Item(452,19).Value = 1071
$ws.Cells.Item(452,20).Value = 7

# Row 453
$ws.Cells.Item(453,1).Value = 11
$ws.Cells.Item(453,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(453,3).Value = "Bíobío"
$ws.Cells.Item(453,4).Value = 44939
$ws.Cells.Item(453,5).Value = 8
$ws.Cells.Item(453,6).Value = "Fruta"
$ws.Cells.Item(453,7).Value = 100101
$ws.Cells.Item(453,8).Value = "Berries"
$ws.Cells.Item(453,9).Value = 100112025
$ws.Cells.Item(453,10).Value = "Frutilla"
$ws.Cells.Item(453,11).Value = "Sin especificar"
$ws.Cells.Item(453,12).Value = "Primera"
$ws.Cells.Item(453,13).Value = 250
$ws.Cells.Item(453,14).Value = 7000
$ws.Cells.Item(453,15).Value = 7000
$ws.Cells.Item(453,16).Value = 7000
$ws.Cells.Item(453,17).Value = "$/caja 7 kilos"
$ws.Cells.Item(453,18).Value = "Región del Maule"
$ws.Cells.Item(453,19).Value = 1000
$ws.Cells.Item(453,20).Value = 7

# Row 454
$ws.Cells.Item(454,1).Value = 11
$ws.Cells.Item(454,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(454,3).Value = "Bíobío"
$ws.Cells.Item(454,4).Value = 44939
$ws.Cells.Item(454,5).Value = 8
$ws.Cells.Item(454,6).Value = "Fruta"
$ws.Cells.Item(454,7).Value = 100101
$ws.Cells.Item(454,8).Value = "Berries"
$ws.Cells.Item(454,9).Value = 100112025
$ws.Cells.Item(454,10).Value = "Frutilla"
$ws.Cells.Item(454,11).Value = "Sin especificar"
$ws.Cells.Item(454,12).Value = "Segunda"
$ws.Cells.Item(454,13).Value = 200
$ws.Cells.Item(454,14).Value = 6000
$ws.Cells.Item(454,15).Value = 6000
$ws.Cells.Item(454,16).Value = 6000
$ws.Cells.Item(454,17).Value = "$/caja 7 kilos"
$ws.Cells.Item(454,18).Value = "Región del Maule"
$ws.Cells.Item(454,19).Value = 857
$ws.Cells.Item(454,20).Value = 7
